# Refresh marketboard-derived profit columns (H,I,J,K,L,M,N) on the
# Yojimbo_Profits leve-crafting sheets (one per DoH job: ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) with the latest scheduled-runner price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Roof Tile
$ws.Range("H19").Value = 585.7143
$ws.Range("I19").Value = 600
$ws.Range("K19").Value = 600
$ws.Range("M19").Value = -425

# Row 125: Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 34578.668
$ws.Range("J125").Value = 1868
$ws.Range("L125").Value = 16812
$ws.Range("N125").Value = -21732

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 2217.7334
$ws.Range("I137").Value = 2428.1538
$ws.Range("J137").Value = 850
$ws.Range("K137").Value = 7284.4614
$ws.Range("L137").Value = 2550
$ws.Range("M137").Value = -4734.4614
$ws.Range("N137").Value = -7650

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Mythril Ingot
$ws.Range("H45").Value = 1567.75
$ws.Range("I45").Value = 1580.5385
$ws.Range("J45").Value = 1544
$ws.Range("K45").Value = 1580.5385
$ws.Range("L45").Value = 1544
$ws.Range("M45").Value = -1203.5385
$ws.Range("N45").Value = -2298

# Row 81: Titanium Headgear of Scouting
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996

# Row 84: Titanium Headgear of Scouting
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

# Row 110: Deepgold Ingot
$ws.Range("H110").Value = 873.6
$ws.Range("I110").Value = 870
$ws.Range("J110").Value = 888
$ws.Range("K110").Value = 870
$ws.Range("L110").Value = 888
$ws.Range("M110").Value = 1175
$ws.Range("N110").Value = -4978

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Iron Ingot
$ws.Range("H20").Value = 1544.4642
$ws.Range("I20").Value = 1431.1818
$ws.Range("K20").Value = 1431.1818
$ws.Range("M20").Value = -1184.1818

# Row 97: High Steel File
$ws.Range("H97").Value = 1485
$ws.Range("I97").Value = 1485
$ws.Range("K97").Value = 1485
$ws.Range("M97").Value = -494

# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 62503110
$ws.Range("I105").Value = 111114930
$ws.Range("J105").Value = 2205.1428
$ws.Range("K105").Value = 111114930
$ws.Range("L105").Value = 2205.1428
$ws.Range("M105").Value = -111113183
$ws.Range("N105").Value = -5699.1428

$ws = $wb.Worksheets.Item("CRP")
# Row 94: Beech Lumber
$ws.Range("H94").Value = 446825.22
$ws.Range("I94").Value = 504453.5
$ws.Range("J94").Value = 400722.6
$ws.Range("K94").Value = 504453.5
$ws.Range("L94").Value = 400722.6
$ws.Range("M94").Value = -504002.5
$ws.Range("N94").Value = -401624.6

# Row 99: Pine Lumber
$ws.Range("H99").Value = 1661.9375
$ws.Range("I99").Value = 1417.75
$ws.Range("J99").Value = 2394.5
$ws.Range("K99").Value = 1417.75
$ws.Range("L99").Value = 2394.5
$ws.Range("M99").Value = 80.25
$ws.Range("N99").Value = -5390.5

# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 1661.9375
$ws.Range("I126").Value = 1417.75
$ws.Range("J126").Value = 2394.5
$ws.Range("K126").Value = 4253.25
$ws.Range("L126").Value = 7183.5
$ws.Range("M126").Value = -1783.25
$ws.Range("N126").Value = -12123.5

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Lavender Oil
$ws.Range("H23").Value = 240.7
$ws.Range("I23").Value = 350
$ws.Range("J23").Value = 167.83333
$ws.Range("K23").Value = 1050
$ws.Range("L23").Value = 503.49999
$ws.Range("M23").Value = -815
$ws.Range("N23").Value = -973.49999

# Row 121: Coffee Biscuit
$ws.Range("H121").Value = 3624117
$ws.Range("I121").Value = 20833662
$ws.Range("J121").Value = 1054.7368
$ws.Range("K121").Value = 62500986
$ws.Range("L121").Value = 3164.2104
$ws.Range("M121").Value = -62499676
$ws.Range("N121").Value = -5784.2104

# Row 125: Borscht
$ws.Range("H125").Value = 5498.3335
$ws.Range("I125").Value = 1995
$ws.Range("J125").Value = 7250
$ws.Range("K125").Value = 5985
$ws.Range("L125").Value = 21750
$ws.Range("M125").Value = -1065
$ws.Range("N125").Value = -31590

# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 1221808.4
$ws.Range("I131").Value = 289.85715
$ws.Range("J131").Value = 1374498.2
$ws.Range("K131").Value = 869.5714499999999
$ws.Range("L131").Value = 4123494.6
$ws.Range("M131").Value = 4170.428550000001
$ws.Range("N131").Value = -4133574.6

$ws = $wb.Worksheets.Item("GSM")
# Row 35: Horn Necklace
$ws.Range("H35").Value = 70000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -70596

# Row 69: Mythrite Needle
$ws.Range("H69").Value = 17687.5
$ws.Range("I69").Value = 8875
$ws.Range("J69").Value = 26500
$ws.Range("K69").Value = 8875
$ws.Range("L69").Value = 26500
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -27998

# Row 72: Mythrite Needle
$ws.Range("H72").Value = 17687.5
$ws.Range("I72").Value = 8875
$ws.Range("J72").Value = 26500
$ws.Range("K72").Value = 26625
$ws.Range("L72").Value = 79500
$ws.Range("M72").Value = -22881
$ws.Range("N72").Value = -86988

# Row 82: Hardsilver Planisphere
$ws.Range("H82").Value = 21439.111
$ws.Range("J82").Value = 21439.111
$ws.Range("L82").Value = 21439.111
$ws.Range("N82").Value = -22205.111

# Row 85: Hardsilver Planisphere
$ws.Range("H85").Value = 21439.111
$ws.Range("J85").Value = 21439.111
$ws.Range("L85").Value = 21439.111
$ws.Range("N85").Value = -24091.111

# Row 102: Durium Ingot
$ws.Range("H102").Value = 1947.8889
$ws.Range("I102").Value = 1916.7142
$ws.Range("J102").Value = 2057
$ws.Range("K102").Value = 1916.7142
$ws.Range("L102").Value = 2057
$ws.Range("M102").Value = -294.7141999999999
$ws.Range("N102").Value = -5301

# Row 112: Diaspore Bracelet of Slaying
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 122: Ametrine
$ws.Range("H122").Value = 967.0909
$ws.Range("I122").Value = 967.0909
$ws.Range("K122").Value = 2901.2727
$ws.Range("M122").Value = -451.2727

# Row 126: Phrygian Gold Ingot
$ws.Range("H126").Value = 3466.6667
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -28940

# Row 132: Lar Ingot
$ws.Range("H132").Value = 12124.857
$ws.Range("I132").Value = 8674.866
$ws.Range("J132").Value = 20749.834
$ws.Range("K132").Value = 26024.598
$ws.Range("L132").Value = 62249.50199999999
$ws.Range("M132").Value = -23494.598
$ws.Range("N132").Value = -67309.50199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 2406.5881
$ws.Range("I7").Value = 2040.4
$ws.Range("J7").Value = 2929.7144
$ws.Range("K7").Value = 2040.4
$ws.Range("L7").Value = 2929.7144
$ws.Range("M7").Value = -1928.4
$ws.Range("N7").Value = -3153.7144

# Row 40: Toad Leather
$ws.Range("H40").Value = 4335.706
$ws.Range("I40").Value = 3654
$ws.Range("J40").Value = 6551.25
$ws.Range("K40").Value = 3654
$ws.Range("L40").Value = 6551.25
$ws.Range("M40").Value = -3518
$ws.Range("N40").Value = -6823.25

# Row 122: Gaja Leather
$ws.Range("H122").Value = 2524.077
$ws.Range("I122").Value = 2182.8
$ws.Range("J122").Value = 3661.6667
$ws.Range("K122").Value = 6548.400000000001
$ws.Range("L122").Value = 10985.0001
$ws.Range("M122").Value = -4098.400000000001
$ws.Range("N122").Value = -15885.0001

# Row 126: Saiga Leather
$ws.Range("H126").Value = 2406.5881
$ws.Range("I126").Value = 2040.4
$ws.Range("J126").Value = 2929.7144
$ws.Range("K126").Value = 6121.200000000001
$ws.Range("L126").Value = 8789.143199999999
$ws.Range("M126").Value = -3651.200000000001
$ws.Range("N126").Value = -13729.1432

# Row 133: Loboskin Amulet of Fending
$ws.Range("H133").Value = 35042.855
$ws.Range("J133").Value = 35042.855
$ws.Range("L133").Value = 35042.855
$ws.Range("N133").Value = -40102.855

# Row 136: Br'aax Leather
$ws.Range("H136").Value = 2888.9795
$ws.Range("I136").Value = 1920.8695
$ws.Range("J136").Value = 3745.3845
$ws.Range("K136").Value = 5762.6085
$ws.Range("L136").Value = 11236.1535
$ws.Range("M136").Value = -3212.6085
$ws.Range("N136").Value = -16336.1535

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Bright Linen Yarn
$ws.Range("H107").Value = 329.8889
$ws.Range("I107").Value = 331.2857
$ws.Range("J107").Value = 325
$ws.Range("K107").Value = 993.8571000000001
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = 926.1428999999999
$ws.Range("N107").Value = -4815

Write-Host "Updated profit columns across $($wb.Worksheets.Count) sheets."
